# Auto-generated-assisted PowerShell COM-interop script
# Applies numeric updates to columns H..N across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
# per the Ixion_Profits.xlsx market-data refresh diff.

$wb = $excel.ActiveWorkbook

# ALC sheet, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1439.5714  # H40: 1891.909 -> 1439.5714
$ws.Cells.Item(40, 9).Value = 1451  # I40: 1928.8572 -> 1451
$ws.Cells.Item(40, 10).Value = 1427  # J40: 1748.2222 -> 1427
$ws.Cells.Item(40, 11).Value = 1451  # K40: 1928.8572 -> 1451
$ws.Cells.Item(40, 12).Value = 1427  # L40: 1748.2222 -> 1427
$ws.Cells.Item(40, 13).Value = -1276  # M40: -1753.8572 -> -1276
$ws.Cells.Item(40, 14).Value = -1777  # N40: -2098.2222 -> -1777

# ALC sheet, row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 1448  # H98: 1423.1111 -> 1448
$ws.Cells.Item(98, 10).Value = 875.25  # J98: 900.2 -> 875.25
$ws.Cells.Item(98, 12).Value = 875.25  # L98: 900.2 -> 875.25
$ws.Cells.Item(98, 14).Value = -3871.25  # N98: -3896.2 -> -3871.25

# ALC sheet, row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 1448  # H122: 1423.1111 -> 1448
$ws.Cells.Item(122, 10).Value = 875.25  # J122: 900.2 -> 875.25
$ws.Cells.Item(122, 12).Value = 2625.75  # L122: 2700.6 -> 2625.75
$ws.Cells.Item(122, 14).Value = -7525.75  # N122: -7600.6 -> -7525.75

# ALC sheet, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1180.2051  # H132: 1385.6774 -> 1180.2051
$ws.Cells.Item(132, 9).Value = 1184.9474  # I132: 1385.6774 -> 1184.9474
$ws.Cells.Item(132, 10).Value = 1000  # J132: 0 -> 1000
$ws.Cells.Item(132, 11).Value = 3554.8422  # K132: 4157.0322 -> 3554.8422
$ws.Cells.Item(132, 12).Value = 3000  # L132: 0 -> 3000
$ws.Cells.Item(132, 13).Value = -1024.8422  # M132: -1627.0322 -> -1024.8422
$ws.Cells.Item(132, 14).Value = -8060  # N132: None -> -8060

# ALC sheet, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2897.8813  # H138: 2480.1487 -> 2897.8813
$ws.Cells.Item(138, 9).Value = 1314.0968  # I138: 1097.0698 -> 1314.0968
$ws.Cells.Item(138, 10).Value = 4651.357  # J138: 4398.613 -> 4651.357
$ws.Cells.Item(138, 11).Value = 3942.2904  # K138: 3291.2094 -> 3942.2904
$ws.Cells.Item(138, 12).Value = 13954.071  # L138: 13195.839 -> 13954.071
$ws.Cells.Item(138, 13).Value = 1197.7096  # M138: 1848.7906 -> 1197.7096
$ws.Cells.Item(138, 14).Value = -24234.071  # N138: -23475.839 -> -24234.071

# ARM sheet, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4131.0728  # H32: 4979.1265 -> 4131.0728
$ws.Cells.Item(32, 9).Value = 3780.8513  # I32: 4851.772 -> 3780.8513
$ws.Cells.Item(32, 11).Value = 3780.8513  # K32: 4851.772 -> 3780.8513
$ws.Cells.Item(32, 13).Value = -3493.8513  # M32: -4564.772 -> -3493.8513

# ARM sheet, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1398.62  # H74: 1633.1025 -> 1398.62
$ws.Cells.Item(74, 9).Value = 1127.7667  # I74: 1327 -> 1127.7667
$ws.Cells.Item(74, 10).Value = 1804.9  # J74: 2073.125 -> 1804.9
$ws.Cells.Item(74, 11).Value = 1127.7667  # K74: 1327 -> 1127.7667
$ws.Cells.Item(74, 12).Value = 1804.9  # L74: 2073.125 -> 1804.9
$ws.Cells.Item(74, 13).Value = -253.7666999999999  # M74: -453 -> -253.7666999999999
$ws.Cells.Item(74, 14).Value = -3552.9  # N74: -3821.125 -> -3552.9

# ARM sheet, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 1398.62  # H77: 1633.1025 -> 1398.62
$ws.Cells.Item(77, 9).Value = 1127.7667  # I77: 1327 -> 1127.7667
$ws.Cells.Item(77, 10).Value = 1804.9  # J77: 2073.125 -> 1804.9
$ws.Cells.Item(77, 11).Value = 5638.8335  # K77: 6635 -> 5638.8335
$ws.Cells.Item(77, 12).Value = 9024.5  # L77: 10365.625 -> 9024.5
$ws.Cells.Item(77, 13).Value = -1270.8335  # M77: -2267 -> -1270.8335
$ws.Cells.Item(77, 14).Value = -17760.5  # N77: -19101.625 -> -17760.5

# ARM sheet, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 3668836.8  # H122: 1605881 -> 3668836.8
$ws.Cells.Item(122, 9).Value = 4279809.5  # I122: 1712739.8 -> 4279809.5
$ws.Cells.Item(122, 11).Value = 12839428.5  # K122: 5138219.4 -> 12839428.5
$ws.Cells.Item(122, 13).Value = -12836978.5  # M122: -5135769.4 -> -12836978.5

# ARM sheet, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 4720.4116  # H132: 7795.6 -> 4720.4116
$ws.Cells.Item(132, 9).Value = 4014  # I132: 10000 -> 4014
$ws.Cells.Item(132, 10).Value = 6415.8  # J132: 7244.5 -> 6415.8
$ws.Cells.Item(132, 11).Value = 12042  # K132: 30000 -> 12042
$ws.Cells.Item(132, 12).Value = 19247.4  # L132: 21733.5 -> 19247.4
$ws.Cells.Item(132, 13).Value = -9512  # M132: -27470 -> -9512
$ws.Cells.Item(132, 14).Value = -24307.4  # N132: -26793.5 -> -24307.4

# BSM sheet, row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 21277616  # H20: 21277618 -> 21277616
$ws.Cells.Item(20, 9).Value = 31250870  # I20: 30303886 -> 31250870
$ws.Cells.Item(20, 10).Value = 1342.4  # J20: 1411.3572 -> 1342.4
$ws.Cells.Item(20, 11).Value = 31250870  # K20: 30303886 -> 31250870
$ws.Cells.Item(20, 12).Value = 1342.4  # L20: 1411.3572 -> 1342.4
$ws.Cells.Item(20, 13).Value = -31250623  # M20: -30303639 -> -31250623
$ws.Cells.Item(20, 14).Value = -1836.4  # N20: -1905.3572 -> -1836.4

# BSM sheet, row 46
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(46, 8).Value = 14100  # H46: 0 -> 14100
$ws.Cells.Item(46, 10).Value = 14100  # J46: 0 -> 14100
$ws.Cells.Item(46, 12).Value = 14100  # L46: 0 -> 14100
$ws.Cells.Item(46, 14).Value = -14696  # N46: None -> -14696

# BSM sheet, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1465.8667  # H94: 1546.5769 -> 1465.8667
$ws.Cells.Item(94, 9).Value = 1167.091  # I94: 1257.6666 -> 1167.091
$ws.Cells.Item(94, 10).Value = 2287.5  # J94: 2760 -> 2287.5
$ws.Cells.Item(94, 11).Value = 1167.091  # K94: 1257.6666 -> 1167.091
$ws.Cells.Item(94, 12).Value = 2287.5  # L94: 2760 -> 2287.5
$ws.Cells.Item(94, 13).Value = -716.0909999999999  # M94: -806.6666 -> -716.0909999999999
$ws.Cells.Item(94, 14).Value = -3189.5  # N94: -3662 -> -3189.5

# BSM sheet, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2318.5898  # H134: 7040.4585 -> 2318.5898
$ws.Cells.Item(134, 9).Value = 2010.1305  # I134: 11834.454 -> 2010.1305
$ws.Cells.Item(134, 10).Value = 2762  # J134: 2984 -> 2762
$ws.Cells.Item(134, 11).Value = 6030.3915  # K134: 35503.362 -> 6030.3915
$ws.Cells.Item(134, 12).Value = 8286  # L134: 8952 -> 8286
$ws.Cells.Item(134, 13).Value = -3495.3915  # M134: -32968.362 -> -3495.3915
$ws.Cells.Item(134, 14).Value = -13356  # N134: -14022 -> -13356

# CRP sheet, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4097.1562  # H31: 3723.7087 -> 4097.1562
$ws.Cells.Item(31, 9).Value = 1966.05  # I31: 1888.7142 -> 1966.05
$ws.Cells.Item(31, 10).Value = 5065.841  # J31: 4388.1035 -> 5065.841
$ws.Cells.Item(31, 11).Value = 1966.05  # K31: 1888.7142 -> 1966.05
$ws.Cells.Item(31, 12).Value = 5065.841  # L31: 4388.1035 -> 5065.841
$ws.Cells.Item(31, 13).Value = -1671.05  # M31: -1593.7142 -> -1671.05
$ws.Cells.Item(31, 14).Value = -5655.841  # N31: -4978.1035 -> -5655.841

# CRP sheet, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 4097.1562  # H34: 3723.7087 -> 4097.1562
$ws.Cells.Item(34, 9).Value = 1966.05  # I34: 1888.7142 -> 1966.05
$ws.Cells.Item(34, 10).Value = 5065.841  # J34: 4388.1035 -> 5065.841
$ws.Cells.Item(34, 11).Value = 1966.05  # K34: 1888.7142 -> 1966.05
$ws.Cells.Item(34, 12).Value = 5065.841  # L34: 4388.1035 -> 5065.841
$ws.Cells.Item(34, 13).Value = -1764.05  # M34: -1686.7142 -> -1764.05
$ws.Cells.Item(34, 14).Value = -5469.841  # N34: -4792.1035 -> -5469.841

# CRP sheet, row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 6955247.5  # H99: 7824562.5 -> 6955247.5
$ws.Cells.Item(99, 9).Value = 9657.4  # I99: 11890.5 -> 9657.4
$ws.Cells.Item(99, 11).Value = 9657.4  # K99: 11890.5 -> 9657.4
$ws.Cells.Item(99, 13).Value = -8159.4  # M99: -10392.5 -> -8159.4

# CRP sheet, row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 6955247.5  # H126: 7824562.5 -> 6955247.5
$ws.Cells.Item(126, 9).Value = 9657.4  # I126: 11890.5 -> 9657.4
$ws.Cells.Item(126, 11).Value = 28972.2  # K126: 35671.5 -> 28972.2
$ws.Cells.Item(126, 13).Value = -26502.2  # M126: -33201.5 -> -26502.2

# CRP sheet, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 2080.8635  # H134: 2414.2778 -> 2080.8635
$ws.Cells.Item(134, 9).Value = 2157.1707  # I134: 2539.394 -> 2157.1707
$ws.Cells.Item(134, 11).Value = 6471.5121  # K134: 7618.181999999999 -> 6471.5121
$ws.Cells.Item(134, 13).Value = -3936.5121  # M134: -5083.181999999999 -> -3936.5121

# CUL sheet, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 139481.89  # H5: 153410.92 -> 139481.89
$ws.Cells.Item(5, 9).Value = 11499  # I5: 17147.834 -> 11499
$ws.Cells.Item(5, 10).Value = 172391.77  # J5: 177457.36 -> 172391.77
$ws.Cells.Item(5, 11).Value = 34497  # K5: 51443.50199999999 -> 34497
$ws.Cells.Item(5, 12).Value = 517175.3099999999  # L5: 532372.08 -> 517175.3099999999
$ws.Cells.Item(5, 13).Value = -34385  # M5: -51331.50199999999 -> -34385
$ws.Cells.Item(5, 14).Value = -517399.3099999999  # N5: -532596.08 -> -517399.3099999999

# CUL sheet, row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 2826.25  # H68: 2682.923 -> 2826.25
$ws.Cells.Item(68, 9).Value = 3884.606  # I68: 3927.0303 -> 3884.606
$ws.Cells.Item(68, 10).Value = 1930.7179  # J68: 1770.5778 -> 1930.7179
$ws.Cells.Item(68, 11).Value = 11653.818  # K68: 11781.0909 -> 11653.818
$ws.Cells.Item(68, 12).Value = 5792.153700000001  # L68: 5311.7334 -> 5792.153700000001
$ws.Cells.Item(68, 13).Value = -10842.818  # M68: -10970.0909 -> -10842.818
$ws.Cells.Item(68, 14).Value = -7414.153700000001  # N68: -6933.7334 -> -7414.153700000001

# CUL sheet, row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(71, 8).Value = 2826.25  # H71: 2682.923 -> 2826.25
$ws.Cells.Item(71, 9).Value = 3884.606  # I71: 3927.0303 -> 3884.606
$ws.Cells.Item(71, 10).Value = 1930.7179  # J71: 1770.5778 -> 1930.7179
$ws.Cells.Item(71, 11).Value = 34961.45400000001  # K71: 35343.2727 -> 34961.45400000001
$ws.Cells.Item(71, 12).Value = 17376.4611  # L71: 15935.2002 -> 17376.4611
$ws.Cells.Item(71, 13).Value = -30905.45400000001  # M71: -31287.2727 -> -30905.45400000001
$ws.Cells.Item(71, 14).Value = -25488.4611  # N71: -24047.2002 -> -25488.4611

# CUL sheet, row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 805.5  # H107: 795.725 -> 805.5
$ws.Cells.Item(107, 9).Value = 285.6154  # I107: 264.57693 -> 285.6154
$ws.Cells.Item(107, 10).Value = 1087.1041  # J107: 1051.463 -> 1087.1041
$ws.Cells.Item(107, 11).Value = 856.8462000000001  # K107: 793.7307900000001 -> 856.8462000000001
$ws.Cells.Item(107, 12).Value = 3261.3123  # L107: 3154.389 -> 3261.3123
$ws.Cells.Item(107, 13).Value = 1063.1538  # M107: 1126.26921 -> 1063.1538
$ws.Cells.Item(107, 14).Value = -7101.3123  # N107: -6994.389 -> -7101.3123

# CUL sheet, row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1071896.4  # H113: 984071.1 -> 1071896.4
$ws.Cells.Item(113, 10).Value = 769817.5600000001  # J113: 556098.25 -> 769817.5600000001
$ws.Cells.Item(113, 12).Value = 2309452.68  # L113: 1668294.75 -> 2309452.68
$ws.Cells.Item(113, 14).Value = -2313792.68  # N113: -1672634.75 -> -2313792.68

# CUL sheet, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 139481.89  # H135: 153410.92 -> 139481.89
$ws.Cells.Item(135, 9).Value = 11499  # I135: 17147.834 -> 11499
$ws.Cells.Item(135, 10).Value = 172391.77  # J135: 177457.36 -> 172391.77
$ws.Cells.Item(135, 11).Value = 103491  # K135: 154330.506 -> 103491
$ws.Cells.Item(135, 12).Value = 1551525.93  # L135: 1597116.24 -> 1551525.93
$ws.Cells.Item(135, 13).Value = -100956  # M135: -151795.506 -> -100956
$ws.Cells.Item(135, 14).Value = -1556595.93  # N135: -1602186.24 -> -1556595.93

# CUL sheet, row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 25580.61  # H139: 4043.1135 -> 25580.61
$ws.Cells.Item(139, 9).Value = 44842.32  # I139: 4969.88 -> 44842.32
$ws.Cells.Item(139, 10).Value = 2650  # J139: 2823.6843 -> 2650
$ws.Cells.Item(139, 11).Value = 134526.96  # K139: 14909.64 -> 134526.96
$ws.Cells.Item(139, 12).Value = 7950  # L139: 8471.052899999999 -> 7950
$ws.Cells.Item(139, 13).Value = -129386.96  # M139: -9769.639999999999 -> -129386.96
$ws.Cells.Item(139, 14).Value = -18230  # N139: -18751.0529 -> -18230

# GSM sheet, row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5737.2  # H70: 5866.5747 -> 5737.2
$ws.Cells.Item(70, 9).Value = 5975.9414  # I70: 5937.9414 -> 5975.9414
$ws.Cells.Item(70, 10).Value = 5229.875  # J70: 5679.923 -> 5229.875
$ws.Cells.Item(70, 11).Value = 5975.9414  # K70: 5937.9414 -> 5975.9414
$ws.Cells.Item(70, 12).Value = 5229.875  # L70: 5679.923 -> 5229.875
$ws.Cells.Item(70, 13).Value = -5705.9414  # M70: -5667.9414 -> -5705.9414
$ws.Cells.Item(70, 14).Value = -5769.875  # N70: -6219.923 -> -5769.875

# GSM sheet, row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 5737.2  # H73: 5866.5747 -> 5737.2
$ws.Cells.Item(73, 9).Value = 5975.9414  # I73: 5937.9414 -> 5975.9414
$ws.Cells.Item(73, 10).Value = 5229.875  # J73: 5679.923 -> 5229.875
$ws.Cells.Item(73, 11).Value = 5975.9414  # K73: 5937.9414 -> 5975.9414
$ws.Cells.Item(73, 12).Value = 5229.875  # L73: 5679.923 -> 5229.875
$ws.Cells.Item(73, 13).Value = -5039.9414  # M73: -5001.9414 -> -5039.9414
$ws.Cells.Item(73, 14).Value = -7101.875  # N73: -7551.923 -> -7101.875

# GSM sheet, row 106
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(106, 8).Value = 45800  # H106: 38000 -> 45800
$ws.Cells.Item(106, 10).Value = 45800  # J106: 38000 -> 45800
$ws.Cells.Item(106, 12).Value = 45800  # L106: 38000 -> 45800
$ws.Cells.Item(106, 14).Value = -48324  # N106: -40524 -> -48324

# GSM sheet, row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 1902.9131  # H107: 2311.0527 -> 1902.9131
$ws.Cells.Item(107, 9).Value = 2675  # I107: 10000 -> 2675
$ws.Cells.Item(107, 10).Value = 1740.3684  # J107: 1883.8889 -> 1740.3684
$ws.Cells.Item(107, 11).Value = 2675  # K107: 10000 -> 2675
$ws.Cells.Item(107, 12).Value = 1740.3684  # L107: 1883.8889 -> 1740.3684
$ws.Cells.Item(107, 13).Value = -755  # M107: -8080 -> -755
$ws.Cells.Item(107, 14).Value = -5580.3684  # N107: -5723.8889 -> -5580.3684

# GSM sheet, row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 5759.375  # H126: 8063.25 -> 5759.375
$ws.Cells.Item(126, 9).Value = 7530.0586  # I126: 9447.076999999999 -> 7530.0586
$ws.Cells.Item(126, 10).Value = 1459.1428  # J126: 2066.6667 -> 1459.1428
$ws.Cells.Item(126, 11).Value = 22590.1758  # K126: 28341.231 -> 22590.1758
$ws.Cells.Item(126, 12).Value = 4377.428400000001  # L126: 6200.000100000001 -> 4377.428400000001
$ws.Cells.Item(126, 13).Value = -20120.1758  # M126: -25871.231 -> -20120.1758
$ws.Cells.Item(126, 14).Value = -9317.428400000001  # N126: -11140.0001 -> -9317.428400000001

# LTW sheet, row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3706300.5  # H22: 5557950.5 -> 3706300.5
$ws.Cells.Item(22, 10).Value = 2876  # J22: 2793.3333 -> 2876
$ws.Cells.Item(22, 12).Value = 2876  # L22: 2793.3333 -> 2876
$ws.Cells.Item(22, 14).Value = -3466  # N22: -3383.3333 -> -3466

# LTW sheet, row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 3706300.5  # H27: 5557950.5 -> 3706300.5
$ws.Cells.Item(27, 10).Value = 2876  # J27: 2793.3333 -> 2876
$ws.Cells.Item(27, 12).Value = 2876  # L27: 2793.3333 -> 2876
$ws.Cells.Item(27, 14).Value = -3090  # N27: -3007.3333 -> -3090

# LTW sheet, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 33334128  # H46: 23810708 -> 33334128
$ws.Cells.Item(46, 9).Value = 33334128  # I46: 33334232 -> 33334128
$ws.Cells.Item(46, 10).Value = 0  # J46: 1900 -> 0
$ws.Cells.Item(46, 11).Value = 33334128  # K46: 33334232 -> 33334128
$ws.Cells.Item(46, 12).Value = 0  # L46: 1900 -> 0
$ws.Cells.Item(46, 13).Value = -33333940  # M46: -33334044 -> -33333940
$ws.Cells.Item(46, 14).ClearContents()  # N46: -2276 -> (removed)

# LTW sheet, row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 17857554  # H55: 16667055 -> 17857554
$ws.Cells.Item(55, 9).Value = 381.72726  # I55: 333.84616 -> 381.72726
$ws.Cells.Item(55, 11).Value = 381.72726  # K55: 333.84616 -> 381.72726
$ws.Cells.Item(55, 13).Value = -208.72726  # M55: -160.84616 -> -208.72726

# LTW sheet, row 119
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(119, 8).Value = 0  # H119: 24500 -> 0
$ws.Cells.Item(119, 10).Value = 0  # J119: 24500 -> 0
$ws.Cells.Item(119, 12).Value = 0  # L119: 24500 -> 0
$ws.Cells.Item(119, 14).ClearContents()  # N119: -34176 -> (removed)

# LTW sheet, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3546538  # H122: 4530742.5 -> 3546538
$ws.Cells.Item(122, 9).Value = 4206021.5  # I122: 5957114 -> 4206021.5
$ws.Cells.Item(122, 11).Value = 12618064.5  # K122: 17871342 -> 12618064.5
$ws.Cells.Item(122, 13).Value = -12615614.5  # M122: -17868892 -> -12615614.5

# LTW sheet, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 27088352  # H132: 30964616 -> 27088352
$ws.Cells.Item(132, 9).Value = 39399132  # I132: 43348464 -> 39399132
$ws.Cells.Item(132, 10).Value = 4639.6  # J132: 4998.25 -> 4639.6
$ws.Cells.Item(132, 11).Value = 118197396  # K132: 130045392 -> 118197396
$ws.Cells.Item(132, 12).Value = 13918.8  # L132: 14994.75 -> 13918.8
$ws.Cells.Item(132, 13).Value = -118194866  # M132: -130042862 -> -118194866
$ws.Cells.Item(132, 14).Value = -18978.8  # N132: -20054.75 -> -18978.8

# LTW sheet, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 3532.8772  # H136: 5201.1724 -> 3532.8772
$ws.Cells.Item(136, 9).Value = 1833.8372  # I136: 4071.568 -> 1833.8372
$ws.Cells.Item(136, 11).Value = 5501.5116  # K136: 12214.704 -> 5501.5116
$ws.Cells.Item(136, 13).Value = -2951.5116  # M136: -9664.704000000002 -> -2951.5116

# WVR sheet, row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 48668.668  # H2: 70003 -> 48668.668
$ws.Cells.Item(2, 9).Value = 6000  # I2: 0 -> 6000
$ws.Cells.Item(2, 11).Value = 6000  # K2: 0 -> 6000
$ws.Cells.Item(2, 13).Value = -5888  # M2: None -> -5888

# WVR sheet, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 3561.3809  # H122: 3716.25 -> 3561.3809
$ws.Cells.Item(122, 9).Value = 2788.8948  # I122: 2918.0557 -> 2788.8948
$ws.Cells.Item(122, 11).Value = 8366.6844  # K122: 8754.167099999999 -> 8366.6844
$ws.Cells.Item(122, 13).Value = -5916.6844  # M122: -6304.167099999999 -> -5916.6844

# WVR sheet, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1786.925  # H132: 1941.6765 -> 1786.925
$ws.Cells.Item(132, 9).Value = 1589.862  # I132: 1786.1305 -> 1589.862
$ws.Cells.Item(132, 10).Value = 2306.4546  # J132: 2266.9092 -> 2306.4546
$ws.Cells.Item(132, 11).Value = 4769.586  # K132: 5358.3915 -> 4769.586
$ws.Cells.Item(132, 12).Value = 6919.3638  # L132: 6800.7276 -> 6919.3638
$ws.Cells.Item(132, 13).Value = -2239.586  # M132: -2828.3915 -> -2239.586
$ws.Cells.Item(132, 14).Value = -11979.3638  # N132: -11860.7276 -> -11979.3638

# WVR sheet, row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(135, 8).Value = 39000  # H135: 0 -> 39000
$ws.Cells.Item(135, 10).Value = 39000  # J135: 0 -> 39000
$ws.Cells.Item(135, 12).Value = 39000  # L135: 0 -> 39000
$ws.Cells.Item(135, 14).Value = -49140  # N135: None -> -49140

# WVR sheet, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 3376.4075  # H136: 2940.8438 -> 3376.4075
$ws.Cells.Item(136, 9).Value = 3682.111  # I136: 3009.652 -> 3682.111
$ws.Cells.Item(136, 11).Value = 11046.333  # K136: 9028.956 -> 11046.333
$ws.Cells.Item(136, 13).Value = -8496.332999999999  # M136: -6478.956 -> -8496.332999999999
